$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 73
$ws.Range("A73").Value = 71
$ws.Range("B73").Value = 7646749
$ws.Range("C73").Value = "Australia ALeague"
$ws.Range("D73").Value = "Australia ALeague"
$ws.Range("E73").Value = 45305.23958333334
$ws.Range("F73").Value = "Brisbane Roar"
$ws.Range("G73").Value = "Newcastle Jets"
$ws.Range("H73").Value = 3
$ws.Range("I73").Value = 2
$ws.Range("J73").Value = "H"
$ws.Range("K73").Value = 1.909
$ws.Range("L73").Value = 4
$ws.Range("M73").Value = 3.4
$ws.Range("N73").Value = 2.4
$ws.Range("O73").Value = 4
$ws.Range("P73").Value = 2.6
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = 1.83
$ws.Range("S73").Value = 2.07
$ws.Range("T73").Value = 3.25
$ws.Range("U73").Value = 1.9
$ws.Range("V73").Value = 1.95
$ws.Range("W73").Value = 1.4
$ws.Range("X73").Value = -1
$ws.Range("Y73").Value = -1
$ws.Range("Z73").Value = 0.8300000000000001
$ws.Range("AA73").Value = -1
$ws.Range("AB73").Value = 0.8999999999999999
$ws.Range("AC73").Value = -1

# Row 74
$ws.Range("A74").Value = 72
$ws.Range("B74").Value = 7646750
$ws.Range("C74").Value = "Australia ALeague"
$ws.Range("D74").Value = "Australia ALeague"
$ws.Range("E74").Value = 45305.23958333334
$ws.Range("F74").Value = "Perth Glory"
$ws.Range("G74").Value = "Wellington Phoenix"
$ws.Range("H74").Value = 3
$ws.Range("I74").Value = 4
$ws.Range("J74").Value = "A"
$ws.Range("K74").Value = 2.45
$ws.Range("L74").Value = 3.75
$ws.Range("M74").Value = 2.55
$ws.Range("N74").Value = 3.1
$ws.Range("O74").Value = 3.8
$ws.Range("P74").Value = 2.05
$ws.Range("Q74").Value = 0.25
$ws.Range("R74").Value = 2
$ws.Range("S74").Value = 1.85
$ws.Range("T74").Value = 3
$ws.Range("U74").Value = 1.925
$ws.Range("V74").Value = 1.925
$ws.Range("W74").Value = -1
$ws.Range("X74").Value = -1
$ws.Range("Y74").Value = 1.05
$ws.Range("Z74").Value = -1
$ws.Range("AA74").Value = 0.8500000000000001
$ws.Range("AB74").Value = 0.925
$ws.Range("AC74").Value = -1

# Row 112
$ws.Range("A112").Value = 110
$ws.Range("B112").Value = 7127379
$ws.Range("C112").Value = "Australia ALeague"
$ws.Range("D112").Value = "Australia ALeague"
$ws.Range("E112").Value = 45347.125
$ws.Range("F112").Value = "Melbourne Victory"
$ws.Range("G112").Value = "Central Coast Mariners"
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = "A"
$ws.Range("K112").Value = 1.95
$ws.Range("L112").Value = 3.6
$ws.Range("M112").Value = 3.8
$ws.Range("N112").Value = 1.909
$ws.Range("O112").Value = 3.6
$ws.Range("P112").Value = 4
$ws.Range("Q112").Value = -0.5
$ws.Range("R112").Value = 1.9
$ws.Range("S112").Value = 1.95
$ws.Range("T112").Value = 2.75
$ws.Range("U112").Value = 1.925
$ws.Range("V112").Value = 1.925
$ws.Range("W112").Value = -1
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = 3
$ws.Range("Z112").Value = -1
$ws.Range("AA112").Value = 0.95
$ws.Range("AB112").Value = -1
$ws.Range("AC112").Value = 0.925

# Row 113
$ws.Range("A113").Value = 111
$ws.Range("B113").Value = 7127376
$ws.Range("C113").Value = "Australia ALeague"
$ws.Range("D113").Value = "Australia ALeague"
$ws.Range("E113").Value = 45347.125
$ws.Range("F113").Value = "Newcastle Jets"
$ws.Range("G113").Value = "Macarthur FC"
$ws.Range("H113").Value = 2
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = "D"
$ws.Range("K113").Value = 1.95
$ws.Range("L113").Value = 4
$ws.Range("M113").Value = 3.4
$ws.Range("N113").Value = 1.909
$ws.Range("O113").Value = 4.2
$ws.Range("P113").Value = 3.6
$ws.Range("Q113").Value = -0.5
$ws.Range("R113").Value = 1.89
$ws.Range("S113").Value = 2.01
$ws.Range("T113").Value = 3.5
$ws.Range("U113").Value = 1.95
$ws.Range("V113").Value = 1.9
$ws.Range("W113").Value = -1
$ws.Range("X113").Value = 3.2
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = -1
$ws.Range("AA113").Value = 1.01
$ws.Range("AB113").Value = 0.95
$ws.Range("AC113").Value = -1

# Row 124
$ws.Range("A124").Value = 122
$ws.Range("B124").Value = 7127388
$ws.Range("C124").Value = "Australia ALeague"
$ws.Range("D124").Value = "Australia ALeague"
$ws.Range("E124").Value = 45361.125
$ws.Range("F124").Value = "Sydney FC"
$ws.Range("G124").Value = "Brisbane Roar"
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = "D"
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 5
$ws.Range("N124").Value = 1.533
$ws.Range("O124").Value = 5.25
$ws.Range("P124").Value = 5
$ws.Range("Q124").Value = -1
$ws.Range("R124").Value = 1.8
$ws.Range("S124").Value = 2.05
$ws.Range("T124").Value = 3.5
$ws.Range("U124").Value = 1.925
$ws.Range("V124").Value = 1.925
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = 4.25
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = 1.05
$ws.Range("AB124").Value = -1
$ws.Range("AC124").Value = 0.925

# Row 125
$ws.Range("A125").Value = 123
$ws.Range("B125").Value = 7128012
$ws.Range("C125").Value = "Australia ALeague"
$ws.Range("D125").Value = "Australia ALeague"
$ws.Range("E125").Value = 45361.125
$ws.Range("F125").Value = "Macarthur FC"
$ws.Range("G125").Value = "Central Coast Mariners"
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = "A"
$ws.Range("K125").Value = 2.4
$ws.Range("L125").Value = 3.5
$ws.Range("M125").Value = 2.75
$ws.Range("N125").Value = 3.4
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 2.05
$ws.Range("Q125").Value = 0.25
$ws.Range("R125").Value = 2.025
$ws.Range("S125").Value = 1.825
$ws.Range("T125").Value = 3
$ws.Range("U125").Value = 2.05
$ws.Range("V125").Value = 1.8
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 1.05
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.825
$ws.Range("AB125").Value = 0
$ws.Range("AC125").Value = -0

# Row 153
$ws.Range("A153").Value = 151
$ws.Range("B153").Value = 7127410
$ws.Range("C153").Value = "Australia ALeague"
$ws.Range("D153").Value = "Australia ALeague"
$ws.Range("E153").Value = 45401.28125
$ws.Range("F153").Value = "Newcastle Jets"
$ws.Range("G153").Value = "Wellington Phoenix"
$ws.Range("H153").Value = 1
$ws.Range("I153").Value = 1
$ws.Range("J153").Value = "D"
$ws.Range("K153").Value = 2.8
$ws.Range("L153").Value = 3.4
$ws.Range("M153").Value = 2.45
$ws.Range("N153").Value = 3
$ws.Range("O153").Value = 4
$ws.Range("P153").Value = 2.15
$ws.Range("Q153").Value = 0.25
$ws.Range("R153").Value = 1.925
$ws.Range("S153").Value = 1.925
$ws.Range("T153").Value = 3
$ws.Range("U153").Value = 1.95
$ws.Range("V153").Value = 1.9
$ws.Range("W153").Value = -1
$ws.Range("X153").Value = 3
$ws.Range("Y153").Value = -1
$ws.Range("Z153").Value = 0.4625
$ws.Range("AA153").Value = -0.5
$ws.Range("AB153").Value = -1
$ws.Range("AC153").Value = 0.8999999999999999

# Row 154
$ws.Range("A154").Value = 152
$ws.Range("B154").Value = 8096897
$ws.Range("C154").Value = "Australia ALeague"
$ws.Range("D154").Value = "Australia ALeague"
$ws.Range("E154").Value = 45402.10416666666
$ws.Range("F154").Value = "Western Sydney Wanderers"
$ws.Range("G154").Value = "Melbourne City"
$ws.Range("K154").Value = 3.25
$ws.Range("L154").Value = 3.8
$ws.Range("M154").Value = 2
$ws.Range("N154").Value = 3.3
$ws.Range("O154").Value = 4
$ws.Range("P154").Value = 2
$ws.Range("Q154").Value = 0.5
$ws.Range("R154").Value = 1.87
$ws.Range("S154").Value = 2.03
$ws.Range("T154").Value = 3.25
$ws.Range("U154").Value = 1.95
$ws.Range("V154").Value = 1.9
$ws.Range("W154").Value = 0
$ws.Range("X154").Value = 0
$ws.Range("Y154").Value = 0
$ws.Range("Z154").Value = 0
$ws.Range("AA154").Value = 0

# Row 155
$ws.Range("A155").Value = 153
$ws.Range("B155").Value = 7127411
$ws.Range("C155").Value = "Australia ALeague"
$ws.Range("D155").Value = "Australia ALeague"
$ws.Range("E155").Value = 45402.1875
$ws.Range("F155").Value = "Melbourne Victory"
$ws.Range("G155").Value = "Brisbane Roar"
$ws.Range("K155").Value = 1.65
$ws.Range("L155").Value = 4
$ws.Range("M155").Value = 4.75
$ws.Range("N155").Value = 1.6
$ws.Range("O155").Value = 4.333
$ws.Range("P155").Value = 5
$ws.Range("Q155").Value = -1
$ws.Range("R155").Value = 2.04
$ws.Range("S155").Value = 1.86
$ws.Range("T155").Value = 3.25
$ws.Range("U155").Value = 2.05
$ws.Range("V155").Value = 1.8
$ws.Range("W155").Value = 0
$ws.Range("X155").Value = 0
$ws.Range("Y155").Value = 0
$ws.Range("Z155").Value = 0
$ws.Range("AA155").Value = 0

# Row 156
$ws.Range("A156").Value = 154
$ws.Range("B156").Value = 7127415
$ws.Range("C156").Value = "Australia ALeague"
$ws.Range("D156").Value = "Australia ALeague"
$ws.Range("E156").Value = 45402.28125
$ws.Range("F156").Value = "Macarthur FC"
$ws.Range("G156").Value = "Sydney FC"
$ws.Range("K156").Value = 3.8
$ws.Range("L156").Value = 4.2
$ws.Range("M156").Value = 1.8
$ws.Range("N156").Value = 4.333
$ws.Range("O156").Value = 4.2
$ws.Range("P156").Value = 1.666
$ws.Range("Q156").Value = 0.75
$ws.Range("R156").Value = 2.02
$ws.Range("S156").Value = 1.88
$ws.Range("T156").Value = 3.5
$ws.Range("U156").Value = 1.925
$ws.Range("V156").Value = 1.925
$ws.Range("W156").Value = 0
$ws.Range("X156").Value = 0
$ws.Range("Y156").Value = 0
$ws.Range("Z156").Value = 0
$ws.Range("AA156").Value = 0

# Row 157
$ws.Range("A157").Value = 155
$ws.Range("B157").Value = 7127414
$ws.Range("C157").Value = "Australia ALeague"
$ws.Range("D157").Value = "Australia ALeague"
$ws.Range("E157").Value = 45403.16666666666
$ws.Range("F157").Value = "Perth Glory"
$ws.Range("G157").Value = "Western United FC"
$ws.Range("K157").Value = 2.4
$ws.Range("L157").Value = 3.6
$ws.Range("M157").Value = 2.625
$ws.Range("N157").Value = 2.4
$ws.Range("O157").Value = 3.75
$ws.Range("P157").Value = 2.7
$ws.Range("Q157").Value = 0
$ws.Range("R157").Value = 1.84
$ws.Range("S157").Value = 2.06
$ws.Range("T157").Value = 3.5
$ws.Range("U157").Value = 1.975
$ws.Range("V157").Value = 1.875
$ws.Range("W157").Value = 0
$ws.Range("X157").Value = 0
$ws.Range("Y157").Value = 0
$ws.Range("Z157").Value = 0
$ws.Range("AA157").Value = 0
